$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the achievement point queries (rows 4-8) ---
# Replace the old single generic "count points" query + a separate ">N" condition
# column with a single combined query per row, and drop the old "Erfuellt wenn" value.
$ws.Range("C4").Value = "SELECT mt.fk_user FROM matchtip mt GROUP BY fk_user HAVIN COUNT(mt.userPoints) > 1"
$ws.Range("D4").ClearContents()

$ws.Range("C5").Value = "SELECT mt.fk_user FROM matchtip mt GROUP BY fk_user HAVIN COUNT(mt.userPoints) > 123"
$ws.Range("D5").ClearContents()

$ws.Range("C6").Value = "SELECT mt.fk_user FROM matchtip mt GROUP BY fk_user HAVIN COUNT(mt.userPoints) > 300"
$ws.Range("D6").ClearContents()

$ws.Range("C7").Value = "SELECT mt.fk_user FROM matchtip mt GROUP BY fk_user HAVIN COUNT(mt.userPoints) > 600"
$ws.Range("D7").ClearContents()

$ws.Range("C8").Value = "SELECT mt.fk_user FROM matchtip mt GROUP BY fk_user HAVIN COUNT(mt.userPoints) > 1234"
$ws.Range("D8").ClearContents()

# --- Highlight the achievement-checker SQL/condition columns in red (rows 10-16) ---
$redRows = 10,11,12,13,14,15,16
foreach ($r in $redRows) {
    $ws.Range("C$r").Font.Color = 255
    $ws.Range("D$r").Font.Color = 255
}

# New blank spacer row (17) formatted the same way as the block above it
$ws.Range("C17").Font.Color = 255
$ws.Range("D17").Font.Color = 255

# Continue the red styling on rows 18-20, and add the (empty) D column cells
# that are now part of the formatted block
$ws.Range("C18").Font.Color = 255
$ws.Range("D18").Font.Color = 255
$ws.Range("D18").ClearContents()

$ws.Range("C19").Font.Color = 255
$ws.Range("D19").Font.Color = 255
$ws.Range("D19").ClearContents()

$ws.Range("C20").Font.Color = 255
$ws.Range("D20").Font.Color = 255

# --- New note row with a large-font reminder ---
$ws.Range("C24").Font.Size = 48
$ws.Range("C24").Value = "NEED TO RETURN USER ID LIST"
$ws.Rows(24).RowHeight = 61.5

# --- Restore the cursor/selection like in the final file ---
$ws.Range("C35").Select()
